# Update "想去人数" (number of people interested) figures for several
# events on the "展览" and "全部类型" sheets, reflecting refreshed scrape
# output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2-8
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 14235
$wsExpo.Range("F3").Value = 333
$wsExpo.Range("F4").Value = 687
$wsExpo.Range("F5").Value = 238
$wsExpo.Range("F6").Value = 554
$wsExpo.Range("F7").Value = 1488
$wsExpo.Range("F8").Value = 141

# Sheet "全部类型" (all types) - rows 2,3,4,5,8,9,11
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14235
$wsAll.Range("F3").Value = 333
$wsAll.Range("F4").Value = 687
$wsAll.Range("F5").Value = 238
$wsAll.Range("F8").Value = 554
$wsAll.Range("F9").Value = 1488
$wsAll.Range("F11").Value = 141

$wb.Save()
